# Added feeding information from 8/6 (2017-08-06) to the "algae-added" sheet.
# This mirrors a new row of observed counts / computed feed volumes, plus a
# new note string, appended right after the existing 7/30-8/5 data (rows 2-8),
# ahead of the trailing blank/formatting row that lived at row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algae-added")

# Seed row 9 with row 7's formatting (date style on A, wrap-text style on B)
# so the new cells pick up the same styles already used by the sheet instead
# of manufacturing brand-new style entries.
$ws.Range("A7:G7").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)

# Date + bucket/strain label + the five raw counts for 8/6.
$ws.Range("A9").Value = 42953
$ws.Range("B9").Value = "500 mL Ciso, 250 609, 250 Chagra"
$ws.Range("C9").Value = 252
$ws.Range("D9").Value = 181
$ws.Range("E9").Value = 197
$ws.Range("F9").Value = 214
$ws.Range("G9").Value = 221

# Extend the same computed columns used by rows 2-8 down into row 9.
$ws.Range("H9").Formula = "=AVERAGE(C9:G9)"
$ws.Range("I9").Formula = "=(H9*9)/0.0009"
$ws.Range("J9").Formula = "=15000*65000"
$ws.Range("K9").Formula = "=J9/I9"
$ws.Range("L9").Value = 500
$ws.Range("M9").Formula = "=L9*I9"
$ws.Range("N9").Formula = "=M9/15000"
$ws.Range("O9").Value = "250 mL Ciso, 125 609, 125 Chagra. Extremely dense, most likely undercounting"

# Matches the author's final cursor position after entering the new row.
[void]$ws.Range("O10").Select()
